$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C1 value
$ws.Range("C1").Value = 9126

# Insert a new row above current row 2 (shifts old row 2 down to row 3)
$ws.Rows.Item(2).Insert()

# Fill in the newly inserted row 2 with the new student's data
$ws.Range("A2").Value = "Сергей"
$ws.Range("B2").Value = "Цыкура"
$ws.Range("C2").Value = 2178

# Update the value in what is now row 3 (previously row 2)
$ws.Range("C3").Value = 5947
